$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 24.1139780380614
$ws.Range("C2").Value = 11.45377727834326
$ws.Range("E2").Value = 10.64562346770988
$ws.Range("F2").Value = 46.03821909381547
$ws.Range("G2").Value = 3.723211791460829
$ws.Range("I2").Value = 32.19702939493263
$ws.Range("J2").Value = 9.456156784107938
$ws.Range("L2").Value = 12.37332449012146
$ws.Range("N2").Value = 20.00034599852896

$ws.Range("B3").Value = 23.71814818957691
$ws.Range("C3").Value = 10.93222397925131
$ws.Range("E3").Value = 10.63468244850968
$ws.Range("F3").Value = 45.94691703261083
$ws.Range("G3").Value = 3.727393897843381
$ws.Range("I3").Value = 32.22735143633449
$ws.Range("J3").Value = 9.486395568345515
$ws.Range("L3").Value = 12.36923263810345
$ws.Range("N3").Value = 20.07548667740893

$ws.Range("B4").Value = 23.4786032000553
$ws.Range("C4").Value = 10.60273340670899
$ws.Range("E4").Value = 10.62791673532077
$ws.Range("F4").Value = 45.90278409889217
$ws.Range("G4").Value = 3.730093802706363
$ws.Range("I4").Value = 32.2538129472429
$ws.Range("J4").Value = 9.5058727764683
$ws.Range("L4").Value = 12.36861842807021
$ws.Range("N4").Value = 20.1236442335587

$ws.Range("B5").Value = 23.38199720109952
$ws.Range("C5").Value = 10.46635368399711
$ws.Range("E5").Value = 10.62514705966246
$ws.Range("F5").Value = 45.88780398422311
$ws.Range("G5").Value = 3.731227378044284
$ws.Range("I5").Value = 32.26656075132582
$ws.Range("J5").Value = 9.51403951036073
$ws.Range("L5").Value = 12.36884678850422
$ws.Range("N5").Value = 20.14377853243046

$ws.Range("B6").Value = 23.36602069845487
$ws.Range("C6").Value = 10.44358743993367
$ws.Range("E6").Value = 10.62468638055327
$ws.Range("F6").Value = 45.88549810349035
$ws.Range("G6").Value = 3.731417624994361
$ws.Range("I6").Value = 32.26879592226473
$ws.Range("J6").Value = 9.515409478057943
$ws.Range("L6").Value = 12.36891365265338
$ws.Range("N6").Value = 20.14715264748376

$ws.Range("B7").Value = 23.47729607354493
$ws.Range("C7").Value = 10.60090237527595
$ws.Range("E7").Value = 10.6278794341646
$ws.Range("F7").Value = 45.90256990250485
$ws.Range("G7").Value = 3.730108955326628
$ws.Range("I7").Value = 32.25397692545194
$ws.Range("J7").Value = 9.505981985345844
$ws.Range("L7").Value = 12.36861956812639
$ws.Range("N7").Value = 20.12391370602229

$ws.Range("B8").Value = 23.97684938019882
$ws.Range("C8").Value = 11.27598831122764
$ws.Range("E8").Value = 10.64186060640893
$ws.Range("F8").Value = 46.00426676357151
$ws.Range("G8").Value = 3.724626448651502
$ws.Range("I8").Value = 32.2058520094397
$ws.Range("J8").Value = 9.466394551159755
$ws.Range("L8").Value = 12.37152050127219
$ws.Range("N8").Value = 20.0258360029787

$ws.Range("B9").Value = 24.97841889333905
$ws.Range("C9").Value = 12.5182340349899
$ws.Range("E9").Value = 10.66892374169393
$ws.Range("F9").Value = 46.29796884992869
$ws.Range("G9").Value = 3.714917248070037
$ws.Range("I9").Value = 32.17403259750148
$ws.Range("J9").Value = 9.395957182264537
$ws.Range("L9").Value = 12.39220570829953
$ws.Range("N9").Value = 19.84946880913529

$ws.Range("B10").Value = 25.71995708805679
$ws.Range("C10").Value = 13.37200240309243
$ws.Range("E10").Value = 10.68862895808964
$ws.Range("F10").Value = 46.57055904465074
$ws.Range("G10").Value = 3.708410751759917
$ws.Range("I10").Value = 32.18917489952731
$ws.Range("J10").Value = 9.348548838902056
$ws.Range("L10").Value = 12.4164492717264
$ws.Range("N10").Value = 19.72952030714517

$ws.Range("B11").Value = 26.05695216065572
$ws.Range("C11").Value = 13.74612088775578
$ws.Range("E11").Value = 10.69756300992845
$ws.Range("F11").Value = 46.70669954893779
$ws.Range("G11").Value = 3.705585104104326
$ws.Range("I11").Value = 32.20449189790898
$ws.Range("J11").Value = 9.327915199933766
$ws.Range("L11").Value = 12.4294182295278
$ws.Range("N11").Value = 19.67702138546051

$ws.Range("B12").Value = 26.18438829840646
$ws.Range("C12").Value = 13.8856393477065
$ws.Range("E12").Value = 10.70094244754938
$ws.Range("F12").Value = 46.75997587857814
$ws.Range("G12").Value = 3.704534263360013
$ws.Range("I12").Value = 32.21150763305528
$ws.Range("J12").Value = 9.320235216153803
$ws.Range("L12").Value = 12.43460584095427
$ws.Range("N12").Value = 19.65743693120093

$ws.Range("B13").Value = 26.15695249296643
$ws.Range("C13").Value = 13.8556886997526
$ws.Range("E13").Value = 10.70021477964989
$ws.Range("F13").Value = 46.74842563084878
$ws.Range("G13").Value = 3.704759730091098
$ws.Range("I13").Value = 32.20994256343127
$ws.Range("J13").Value = 9.321883308646894
$ws.Range("L13").Value = 12.43347633808815
$ws.Range("N13").Value = 19.66164166209346

$ws.Range("B14").Value = 26.06744066590744
$ws.Range("C14").Value = 13.75764280642232
$ws.Range("E14").Value = 10.6978411087839
$ws.Range("F14").Value = 46.71104823467183
$ws.Range("G14").Value = 3.705498267312821
$ws.Range("I14").Value = 32.20504470941432
$ws.Range("J14").Value = 9.3272806902895
$ws.Range("L14").Value = 12.42983948992279
$ws.Range("N14").Value = 19.67540424122527

$ws.Range("B15").Value = 26.01258536087925
$ws.Range("C15").Value = 13.69730388576775
$ws.Range("E15").Value = 10.69638670412023
$ws.Range("F15").Value = 46.68837715288691
$ws.Range("G15").Value = 3.705953135553634
$ws.Range("I15").Value = 32.20220301544158
$ws.Range("J15").Value = 9.330604114335166
$ws.Range("L15").Value = 12.42764775311937
$ws.Range("N15").Value = 19.683872692851

$ws.Range("B16").Value = 25.69791584478111
$ws.Range("C16").Value = 13.34725683821348
$ws.Range("E16").Value = 10.68804454632974
$ws.Range("F16").Value = 46.5619043253018
$ws.Range("G16").Value = 3.708598103976157
$ws.Range("I16").Value = 32.18834382501447
$ws.Range("J16").Value = 9.349916012734564
$ws.Range("L16").Value = 12.41564055856959
$ws.Range("N16").Value = 19.73299269538806

$ws.Range("B17").Value = 25.50469771891508
$ws.Range("C17").Value = 13.12878366894115
$ws.Range("E17").Value = 10.68291966538573
$ws.Range("F17").Value = 46.4874110038168
$ws.Range("G17").Value = 3.710254986198379
$ws.Range("I17").Value = 32.18200324419443
$ws.Range("J17").Value = 9.362001691467176
$ws.Range("L17").Value = 12.40876982173402
$ws.Range("N17").Value = 19.76365445829607

$ws.Range("B18").Value = 25.39353969711454
$ws.Range("C18").Value = 13.00178588684921
$ws.Range("E18").Value = 10.67996911323885
$ws.Range("F18").Value = 46.44570841969215
$ws.Range("G18").Value = 3.711220619135374
$ws.Range("I18").Value = 32.17914951233087
$ws.Range("J18").Value = 9.369040874625489
$ws.Range("L18").Value = 12.40500076021094
$ws.Range("N18").Value = 19.7814848602908

$ws.Range("B19").Value = 25.35590336093206
$ws.Range("C19").Value = 12.95856019612439
$ws.Range("E19").Value = 10.67896959790783
$ws.Range("F19").Value = 46.43178576523066
$ws.Range("G19").Value = 3.711549740296408
$ws.Range("I19").Value = 32.17831939847718
$ws.Range("J19").Value = 9.371439322342779
$ws.Range("L19").Value = 12.40375609005607
$ws.Range("N19").Value = 19.78755538725719

$ws.Range("B20").Value = 25.52526955524222
$ws.Range("C20").Value = 13.15217974951671
$ws.Range("E20").Value = 10.68346550889824
$ws.Range("F20").Value = 46.4952226993057
$ws.Range("G20").Value = 3.710077301173868
$ws.Range("I20").Value = 32.18259608802405
$ws.Range("J20").Value = 9.36070606554053
$ws.Range("L20").Value = 12.40948232119346
$ws.Range("N20").Value = 19.76037033438397

$ws.Range("B21").Value = 26.09373825136326
$ws.Range("C21").Value = 13.78650039225917
$ws.Range("E21").Value = 10.69853840878656
$ws.Range("F21").Value = 46.7219803170582
$ws.Range("G21").Value = 3.705280821740884
$ws.Range("I21").Value = 32.20645031614615
$ws.Range("J21").Value = 9.325691729201498
$ws.Range("L21").Value = 12.43090023570977
$ws.Range("N21").Value = 19.6713538224959

$ws.Range("B22").Value = 26.46418014739242
$ws.Range("C22").Value = 14.18848580327218
$ws.Range("E22").Value = 10.70836817953515
$ws.Range("F22").Value = 46.88020936371918
$ws.Range("G22").Value = 3.702257728594821
$ws.Range("I22").Value = 32.2291263262762
$ws.Range("J22").Value = 9.303585870554141
$ws.Range("L22").Value = 12.44650872752574
$ws.Range("N22").Value = 19.6148996043372

$ws.Range("B23").Value = 26.26660878818429
$ws.Range("C23").Value = 13.97511888992869
$ws.Range("E23").Value = 10.70312358481542
$ws.Range("F23").Value = 46.79484997378407
$ws.Range("G23").Value = 3.703861032686165
$ws.Range("I23").Value = 32.21637445963469
$ws.Range("J23").Value = 9.315313184788705
$ws.Range("L23").Value = 12.43803167863331
$ws.Range("N23").Value = 19.64487305908222

$ws.Range("B24").Value = 25.51596924981247
$ws.Range("C24").Value = 13.14160671768611
$ws.Range("E24").Value = 10.6832187461426
$ws.Range("F24").Value = 46.49168752716407
$ws.Range("G24").Value = 3.710157591856181
$ws.Range("I24").Value = 32.1823255980953
$ws.Range("J24").Value = 9.361291534640461
$ws.Range("L24").Value = 12.40915963617549
$ws.Range("N24").Value = 19.76185445581468

$ws.Range("B25").Value = 24.70592026221225
$ws.Range("C25").Value = 12.19194566344512
$ws.Range("E25").Value = 10.66163646601505
$ws.Range("F25").Value = 46.20847870473785
$ws.Range("G25").Value = 3.717433159341276
$ws.Range("I25").Value = 32.17590041872288
$ws.Range("J25").Value = 9.414246757237915
$ws.Range("L25").Value = 12.38501347552746
$ws.Range("N25").Value = 19.89548256542346
